$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (e.g. "330.30", "26.459.63").
# Excel auto-converts plain numeric-looking strings assigned via .Value
# into real numbers (dropping trailing zeros, changing the cell type),
# so for the cells whose new price would be mis-detected as a number we
# pre-format as Text ("@") to keep them as the exact literal string.
$fmtCells = @("D4","D5","D6","D7","D9","D11","D12","D13","D14","D15","D17","D18","D19","D21","D22","D23","D25","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $fmtCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "26.459.63"
$ws.Range("E2").Value = "  +6.24%  "
$ws.Range("D3").Value = "1.721.64"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "330.69"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.3714"
$ws.Range("E7").Value = "  +2.35%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").Value = "0.3348"
$ws.Range("E9").Value = "  +2.47%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").Value = "0.07359"
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "6.358"
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("D14").Value = "20.05"
$ws.Range("E14").Value = "  +2.64%  "
$ws.Range("D15").Value = "7.032"
$ws.Range("E15").Value = "  +6.47%  "
$ws.Range("D16").Value = "1.721.49"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("D17").Value = "0.00001068"
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").Value = "0.06630"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "82.27"
$ws.Range("E19").Value = "  +4.04%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "16.50"
$ws.Range("E21").Value = "  +4.60%  "
$ws.Range("D22").Value = "6.104"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "12.77"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "26.427.99"
$ws.Range("E24").Value = "  +6.28%  "
$ws.Range("D25").Value = "2.440"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  +17.20%  "
$ws.Range("D27").Value = "2.377"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "152.18"
$ws.Range("E28").Value = "  +2.36%  "
$ws.Range("D29").Value = "19.36"
$ws.Range("E29").Value = "  +3.86%  "
$ws.Range("D30").Value = "1.921.46"
$ws.Range("E30").Value = "  +4.62%  "
$ws.Range("D31").Value = "130.68"
$ws.Range("E31").Value = "  +4.33%  "
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  +1.01%  "
$ws.Range("D33").Value = "5.941"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").Value = "0.08578"
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").Value = "1.697"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").Value = "12.64"
$ws.Range("E36").Value = "  +3.50%  "
$ws.Range("D37").Value = "5.340"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("D38").Value = "0.02314"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "0.2149"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").Value = "0.06186"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").Value = "8.418"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").Value = "1.221"
$ws.Range("E42").Value = "  -4.79%  "
$ws.Range("D43").Value = "0.6151"
$ws.Range("E43").Value = "  +3.69%  "
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.28%  "
$ws.Range("D45").Value = "13.99"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("D46").Value = "3.888"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "0.5946"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("D48").Value = "127.48"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "2.028"
$ws.Range("E49").Value = "  +4.05%  "
$ws.Range("D50").Value = "0.07175"
$ws.Range("E50").Value = "  +2.74%  "
$ws.Range("E51").Value = "  +2.08%  "
